$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 160
$ws.Range("I12").Value = 140
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 140
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = -540

$ws.Range("H70").Value = 3205.0588
$ws.Range("I70").Value = 1742.8572
$ws.Range("J70").Value = 4228.6
$ws.Range("K70").Value = 5228.571599999999
$ws.Range("L70").Value = 12685.8
$ws.Range("M70").Value = -4958.571599999999
$ws.Range("N70").Value = -13225.8

$ws.Range("H73").Value = 3205.0588
$ws.Range("I73").Value = 1742.8572
$ws.Range("J73").Value = 4228.6
$ws.Range("K73").Value = 5228.571599999999
$ws.Range("L73").Value = 12685.8
$ws.Range("M73").Value = -4292.571599999999
$ws.Range("N73").Value = -14557.8

$ws.Range("H74").Value = 12500
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14064

$ws.Range("H77").Value = 12500
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 75000
$ws.Range("M77").Value = -70320

$ws.Range("H113").Value = 7382.722
$ws.Range("J113").Value = 7998.636
$ws.Range("L113").Value = 7998.636
$ws.Range("N113").Value = -14506.636

$ws.Range("H132").Value = 10330.1875
$ws.Range("I132").Value = 1761.537
$ws.Range("K132").Value = 5284.611
$ws.Range("M132").Value = -2754.611

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14728.167
$ws.Range("I2").Value = 17284.3
$ws.Range("K2").Value = 17284.3
$ws.Range("M2").Value = -17171.3

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H63").Value = 11625
$ws.Range("I63").Value = 6600
$ws.Range("K63").Value = 6600
$ws.Range("M63").Value = -5914

$ws.Range("H66").Value = 11625
$ws.Range("I66").Value = 6600
$ws.Range("K66").Value = 33000
$ws.Range("M66").Value = -29568

$ws.Range("H74").Value = 3723.682
$ws.Range("I74").Value = 4237.6875
$ws.Range("K74").Value = 4237.6875
$ws.Range("M74").Value = -3363.6875

$ws.Range("H77").Value = 3723.682
$ws.Range("I77").Value = 4237.6875
$ws.Range("K77").Value = 21188.4375
$ws.Range("M77").Value = -16820.4375

$ws.Range("H116").Value = 14728.167
$ws.Range("I116").Value = 17284.3
$ws.Range("K116").Value = 17284.3
$ws.Range("M116").Value = -14990.3

$ws.Range("H122").Value = 1288.625
$ws.Range("I122").Value = 1288.625
$ws.Range("K122").Value = 3865.875
$ws.Range("M122").Value = -1415.875

$ws.Range("H132").Value = 3527.7896
$ws.Range("I132").Value = 3502.5833
$ws.Range("K132").Value = 10507.7499
$ws.Range("M132").Value = -7977.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14728.167
$ws.Range("I3").Value = 17284.3
$ws.Range("K3").Value = 17284.3
$ws.Range("M3").Value = -17170.3

$ws.Range("H20").Value = 1357.4
$ws.Range("I20").Value = 1410.3846
$ws.Range("J20").Value = 1259
$ws.Range("K20").Value = 1410.3846
$ws.Range("L20").Value = 1259
$ws.Range("M20").Value = -1163.3846
$ws.Range("N20").Value = -1753

$ws.Range("H105").Value = 3707.3333
$ws.Range("I105").Value = 1773.7
$ws.Range("J105").Value = 6124.375
$ws.Range("K105").Value = 1773.7
$ws.Range("L105").Value = 6124.375
$ws.Range("M105").Value = -26.70000000000005
$ws.Range("N105").Value = -9618.375

$ws.Range("H134").Value = 1888.7736
$ws.Range("I134").Value = 1319.907
$ws.Range("J134").Value = 4334.9
$ws.Range("K134").Value = 3959.721
$ws.Range("L134").Value = 13004.7
$ws.Range("M134").Value = -1424.721
$ws.Range("N134").Value = -18074.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7823.25
$ws.Range("I62").Value = 8834.333000000001
$ws.Range("J62").Value = 4790
$ws.Range("K62").Value = 8834.333000000001
$ws.Range("L62").Value = 4790
$ws.Range("M62").Value = -8210.333000000001
$ws.Range("N62").Value = -6038

$ws.Range("H65").Value = 7823.25
$ws.Range("I65").Value = 8834.333000000001
$ws.Range("J65").Value = 4790
$ws.Range("K65").Value = 44171.665
$ws.Range("L65").Value = 23950
$ws.Range("M65").Value = -41051.665
$ws.Range("N65").Value = -30190

$ws.Range("H99").Value = 24937918
$ws.Range("I99").Value = 8131918.5
$ws.Range("J99").Value = 33340918
$ws.Range("K99").Value = 8131918.5
$ws.Range("L99").Value = 33340918
$ws.Range("M99").Value = -8130420.5
$ws.Range("N99").Value = -33343914

$ws.Range("H126").Value = 24937918
$ws.Range("I126").Value = 8131918.5
$ws.Range("J126").Value = 33340918
$ws.Range("K126").Value = 24395755.5
$ws.Range("L126").Value = 100022754
$ws.Range("M126").Value = -24393285.5
$ws.Range("N126").Value = -100027694

$ws.Range("H132").Value = 2239.1428
$ws.Range("I132").Value = 2179
$ws.Range("K132").Value = 6537
$ws.Range("M132").Value = -4007

$ws.Range("H134").Value = 2110.1333
$ws.Range("I134").Value = 1776
$ws.Range("J134").Value = 3029
$ws.Range("K134").Value = 5328
$ws.Range("L134").Value = 9087
$ws.Range("M134").Value = -2793
$ws.Range("N134").Value = -14157

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2333.3333
$ws.Range("I121").Value = 1500
$ws.Range("J121").Value = 2750
$ws.Range("K121").Value = 4500
$ws.Range("L121").Value = 8250
$ws.Range("M121").Value = -3190
$ws.Range("N121").Value = -10870

$ws.Range("H129").Value = 1378.125
$ws.Range("I129").Value = 489.44446
$ws.Range("J129").Value = 2520.7144
$ws.Range("K129").Value = 1468.33338
$ws.Range("L129").Value = 7562.1432
$ws.Range("M129").Value = 3531.66662
$ws.Range("N129").Value = -17562.1432

$ws.Range("H131").Value = 4903.8
$ws.Range("I131").Value = 3264.5
$ws.Range("K131").Value = 9793.5
$ws.Range("M131").Value = -4753.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws.Range("H141").Value = 83948
$ws.Range("J141").Value = 83948
$ws.Range("L141").Value = 83948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2735.158
$ws.Range("I55").Value = 2985.625
$ws.Range("J55").Value = 2553
$ws.Range("K55").Value = 2985.625
$ws.Range("L55").Value = 2553
$ws.Range("M55").Value = -2812.625
$ws.Range("N55").Value = -2899

$ws.Range("H93").Value = 2741.6667
$ws.Range("I93").Value = 3144.4443
$ws.Range("K93").Value = 3144.4443
$ws.Range("M93").Value = -1896.4443

$ws.Range("H122").Value = 5907.5
$ws.Range("I122").Value = 2720
$ws.Range("J122").Value = 6704.375
$ws.Range("K122").Value = 8160
$ws.Range("L122").Value = 20113.125
$ws.Range("M122").Value = -5710
$ws.Range("N122").Value = -25013.125

$ws.Range("H132").Value = 3913.973
$ws.Range("I132").Value = 2410.7856
$ws.Range("J132").Value = 8590.556
$ws.Range("K132").Value = 7232.3568
$ws.Range("L132").Value = 25771.668
$ws.Range("M132").Value = -4702.3568
$ws.Range("N132").Value = -30831.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 710.5599999999999
$ws.Range("I113").Value = 424.8421
$ws.Range("J113").Value = 1615.3334
$ws.Range("K113").Value = 1274.5263
$ws.Range("L113").Value = 4846.0002
$ws.Range("M113").Value = 895.4737
$ws.Range("N113").Value = -9186.0002

$ws.Range("H132").Value = 3077.8484
$ws.Range("I132").Value = 2376.2593
$ws.Range("J132").Value = 6235
$ws.Range("K132").Value = 7128.777900000001
$ws.Range("L132").Value = 18705
$ws.Range("M132").Value = -4598.777900000001
$ws.Range("N132").Value = -23765

$ws.Range("H140").Value = 72547.25
$ws.Range("I140").Value = 48090
$ws.Range("J140").Value = 80699.664
$ws.Range("K140").Value = 48090
$ws.Range("L140").Value = 80699.664
$ws.Range("M140").Value = -42910
$ws.Range("N140").Value = -91059.664

$ws.Range("H141").Value = 97434.266
$ws.Range("J141").Value = 97965.28999999999
$ws.Range("L141").Value = 97965.28999999999
$ws.Range("N141").Value = -108325.29
